$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

# Fill in the fuzzy_4 ("4_fuzzy") test results on row 6 (columns G:J and M:P)
# to match the values/style already used by F6/K6 and L6.
$ws.Range("G6:J6").Value = 1
$ws.Range("M6:P6").Value = 0.25

$ws.Range("F6").Copy()
$ws.Range("G6:J6").PasteSpecial(-4122)

$ws.Range("L6").Copy()
$ws.Range("M6:P6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update the active sheet's view/selection: drop the frozen/topLeft cell
# override and move the selection from L6 to P6.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P6").Select()
